$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) from 45653 to 45654 for rows 2 through 37
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45653) {
        $cell.Value2 = 45654
    }
}

# Row 36 and 37: swap A (Beteckning) and G (Area) values
$ws.Range("A36").Value = "A 60500-2024"
$ws.Range("G36").Value = 0.8

$ws.Range("A37").Value = "A 60501-2024"
$ws.Range("G37").Value = 0.6
